$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (player, position(s), team) for rows 2-19.
$data = New-Object 'object[,]' 18,3

$data[0,0]  = "Scoot Henderson";        $data[0,1]  = "PG";         $data[0,2]  = "Portland Trail Blazers"
$data[1,0]  = "James Harden";           $data[1,1]  = "PG,SG";      $data[1,2]  = "LA Clippers"
$data[2,0]  = "Anfernee Simons";        $data[2,1]  = "PG,SG";      $data[2,2]  = "Portland Trail Blazers"
$data[3,0]  = "Fred VanVleet";          $data[3,1]  = "PG";         $data[3,2]  = "Houston Rockets"
$data[4,0]  = "Anthony Edwards";        $data[4,1]  = "SG,SF";      $data[4,2]  = "Minnesota Timberwolves"
$data[5,0]  = "Paul George";            $data[5,1]  = "SG,SF,PF";   $data[5,2]  = "Philadelphia 76ers"
$data[6,0]  = "Giannis Antetokounmpo";  $data[6,1]  = "PF,C";       $data[6,2]  = "Milwaukee Bucks"
$data[7,0]  = "Grayson Allen";          $data[7,1]  = "PG,SG,SF";   $data[7,2]  = "Phoenix Suns"
$data[8,0]  = "Jaren Jackson Jr.";      $data[8,1]  = "PF,C";       $data[8,2]  = "Memphis Grizzlies"
$data[9,0]  = "Zach Edey";              $data[9,1]  = "C";          $data[9,2]  = "Memphis Grizzlies"
$data[10,0] = "Ivica Zubac";            $data[10,1] = "C";          $data[10,2] = "LA Clippers"
$data[11,0] = "Andrew Wiggins";         $data[11,1] = "SF,PF";      $data[11,2] = "Golden State Warriors"
$data[12,0] = "Jayson Tatum";           $data[12,1] = "SF,PF";      $data[12,2] = "Boston Celtics"
$data[13,0] = "Dennis Schröder";        $data[13,1] = "PG,SG";      $data[13,2] = "Golden State Warriors"
$data[14,0] = "Amen Thompson";          $data[14,1] = "SG,SF";      $data[14,2] = "Houston Rockets"
$data[15,0] = "Zion Williamson";        $data[15,1] = "PF,C";       $data[15,2] = "New Orleans Pelicans"
$data[16,0] = "Jonathan Kuminga";       $data[16,1] = "SF,PF";      $data[16,2] = "Golden State Warriors"
$data[17,0] = "Jerami Grant";           $data[17,1] = "SF,PF";      $data[17,2] = "Portland Trail Blazers"

$ws.Range("A2:C19").Value = $data
